$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# --- New "Sprint 6 backlog" section (rows 58-63) -----------------------

# Section title
$ws.Cells.Item(58, 2).Value = "Sprint 6 backlog"

# Header row (same style as the other sprint tables: centered, bordered -> style index 3)
$headerSrc = $ws.Range("B56:H56")
$headerDst = $ws.Range("B59:H59")
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(59, 2).Value = "Task Number"
$ws.Cells.Item(59, 3).Value = "Task Name"
$ws.Cells.Item(59, 4).Value = "Task Owner"
$ws.Cells.Item(59, 5).Value = "Task Reviewer"
$ws.Cells.Item(59, 6).Value = "To-Do"
$ws.Cells.Item(59, 7).Value = "Reviewing"
$ws.Cells.Item(59, 8).Value = "Done"

# Task rows 1 & 2 formatting (copy from an existing bordered/centered row)
$taskSrc = $ws.Range("B56:H56")
$row60Dst = $ws.Range("B60:H60")
$taskSrc.Copy()
$row60Dst.PasteSpecial(-4122)
$row61Dst = $ws.Range("B61:H61")
$taskSrc.Copy()
$row61Dst.PasteSpecial(-4122)

# Task numbers
$ws.Cells.Item(60, 2).Value = 1
$ws.Cells.Item(61, 2).Value = 2

# Task names (entered for both rows before task owner, matching original authoring order)
$ws.Cells.Item(60, 3).Value = "Implement Webcam"
$ws.Cells.Item(61, 3).Value = "Test Webcam"

# Task owners
$ws.Cells.Item(60, 4).Value = "António/Daniel"
$ws.Cells.Item(61, 4).Value = "António/Daniel"

# Task reviewers
$ws.Cells.Item(60, 5).Value = "N/A"
$ws.Cells.Item(61, 5).Value = "N/A"

# Status marks: row1 Done (H), row2 Reviewing (F)
$ws.Cells.Item(60, 8).Value = "X"
$ws.Cells.Item(61, 6).Value = "X"

$excel.CutCopyMode = 0

# Trailing blank, centered rows (new style: no border, centered -> style index 4)
$blank = $ws.Range("B62:H63")
$blank.HorizontalAlignment = -4108  # xlCenter

# --- Selection / view state ---------------------------------------------
$ws.Range("K59").Select()
